$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicated "Contact" row (row 11) - rows below shift up by one.
$ws.Rows.Item(11).Delete()

# Update Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Update Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Replace the remaining "Contact" row with Jurisdiction info
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Update Case Sensitive value (force literal text "true", not boolean TRUE)
$ws.Range("B14").Formula = '="true"'
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4163) | Out-Null
